$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "Author: " + "Kenny" (two runs) -> single run "Author: Lion"
# This also removes the old "_GoBack" bookmark that used to sit right
# after the "Kenny" run, since Find/Execute rewrites that stretch of
# the paragraph as one run and the bookmark fell inside the replaced
# range.
# ---------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("Author: Kenny", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Author: Lion", 2)

# ---------------------------------------------------------------------
# Change 2: re-create the "_GoBack" bookmark (collapsed) right after the
# trailing " " run of the "Next ,  " paragraph (the paragraph right
# before the document's final, empty paragraph).
#
# Inserting a collapsed bookmark exactly at a position that is
# immediately before an empty paragraph lands it at the wrong spot, so
# a one-character placeholder is inserted first to give the bookmark a
# safe anchor, then removed again once the bookmark is in place - the
# bookmark (being zero-width) stays put when the placeholder is deleted.
# ---------------------------------------------------------------------
$paraRng = $d.Content
$paraRng.Find.Execute("Next ,", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$targetPara = $paraRng.Paragraphs(1)
$paraEnd = $targetPara.Range.End - 1

$placeholder = $d.Range($paraEnd, $paraEnd)
$placeholder.InsertAfter("X")

$bookmarkRng = $d.Range($paraEnd, $paraEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRng)

$placeholderRng = $d.Range($paraEnd, $paraEnd + 1)
$placeholderRng.Delete()

Write-Host "Done"
